$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Section properties: restart page numbering at 3034  (<w:pgNumType w:start=".."/>)
# ---------------------------------------------------------------------------
$d.Sections(1).Footers.Item(1).PageNumbers.StartingNumber = 3034

# ---------------------------------------------------------------------------
# All edits below are applied from the LAST paragraph to the FIRST so that
# paragraph indices already visited never shift under us. Only the VISTO and
# CONSIDERANDO paragraphs are split into two, and both sit above (i.e. are
# indexed lower than) every other paragraph we touch, so doing those splits
# last keeps every earlier "Paragraphs(N)" lookup valid.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Paragraph 10: "ARTICULO TERCERO: COMUNIQUESE, COPIESE Y ARCHIVESE."
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(10)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0
$base = $p.Range.Start
$d.Range($base, $base + 16).Font.Underline = 1
$d.Range($base + 16, $base + 17).Font.Underline = 1

# ---------------------------------------------------------------------------
# Paragraph 9: "ARTICULO SEGUNDO: El mencionado Articulo dispondrá..."
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(9)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0
$base = $p.Range.Start
$d.Range($base, $base + 16).Font.Underline = 1
$d.Range($base + 16, $base + 17).Font.Underline = 1

# ---------------------------------------------------------------------------
# Paragraph 8: "ARTICULO PRIMERO: AGREGUESEa la Ordenanza Nº 1649 ... Bis."
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(8)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0
$base = $p.Range.Start
$d.Range($base, $base + 16).Font.Underline = 1
$d.Range($base + 16, $base + 17).Font.Underline = 1

# collapse the long run of spaces before "(Reglamento..." down to one space
$full = $p.Range.Text
$idxSp = $full.IndexOf("                    (")
$rSp = $d.Range($base + $idxSp, $base + $idxSp + 21)
$rSp.Text = " ("

# ---------------------------------------------------------------------------
# Paragraph 7: "POR EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA"
#   -> drop the leading "POR " run, bold + center + indent the rest
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(7)
$base = $p.Range.Start
$d.Range($base, $base + 4).Delete()

$p.Format.KeepWithNext = 1
$p.Format.SpaceBefore = 18
$p.Format.SpaceAfter = 18
$p.Format.LeftIndent = 99.2
$p.Format.RightIndent = 99.2
$p.Format.Alignment = 1
$p.Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# Paragraph 6: "Sin embargo muchas veces..." - drop the leading space
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(6)
$base = $p.Range.Start
$d.Range($base, $base + 1).Delete()
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# ---------------------------------------------------------------------------
# Paragraph 5: "Que los pedidos de trabajos..." - drop the leading space
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(5)
$base = $p.Range.Start
$d.Range($base, $base + 1).Delete()
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 6
$p.Format.Alignment = 0

# ---------------------------------------------------------------------------
# Paragraph 4: "CONSIDERANDO: Que los proyectos de Minuta..."
#   -> split into "CONSIDERANDO: " (bold) and the body (new leading space)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(4)
$base = $p.Range.Start
$d.Range($base + 14, $base + 14).InsertParagraphAfter()

$p4 = $d.Paragraphs(4)
$p4.Format.KeepWithNext = 1
$p4.Format.SpaceBefore = 12
$p4.Format.SpaceAfter = 6
$p4.Format.Alignment = 0
$p4.Range.Font.Bold = 1

$p5 = $d.Paragraphs(5)
$p5.Format.KeepWithNext = 1
$p5.Format.SpaceAfter = 6
$p5.Format.Alignment = 0
$d.Range($p5.Range.Start, $p5.Range.Start).InsertBefore(" ")

# ---------------------------------------------------------------------------
# Paragraph 3: "VISTO: Que es necesario darle..."
#   -> split into "VISTO: " (bold) and the body (new leading space)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(3)
$base = $p.Range.Start
$d.Range($base + 7, $base + 7).InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$p3.Format.KeepWithNext = 1
$p3.Format.SpaceBefore = 12
$p3.Format.SpaceAfter = 6
$p3.Format.Alignment = 0
$p3.Range.Font.Bold = 1

$p4b = $d.Paragraphs(4)
$p4b.Format.KeepWithNext = 1
$p4b.Format.SpaceAfter = 6
$p4b.Format.Alignment = 0
$d.Range($p4b.Range.Start, $p4b.Range.Start).InsertBefore(" ")

# ---------------------------------------------------------------------------
# Paragraph 2: "ORDENANZA Nº 2040"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(2)
$p.Format.KeepWithNext = 1
$p.Format.SpaceBefore = 12
$p.Format.SpaceAfter = 18
$p.Range.Font.Bold = 1

# ---------------------------------------------------------------------------
# Paragraph 1: "Yerba Buena, 07 de Abril de 2016"
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(1)
$p.Format.KeepWithNext = 1
$p.Format.SpaceAfter = 12
